# Code v1.0 updated for gratuity calculator test
# Update the "Gold Price" sheet with refreshed price figures.
# Values are stored as literal text (e.g. "48501.00", "-0.82%"), so a
# leading apostrophe is used to force text entry and keep Excel from
# re-interpreting these look-like-numbers strings as numeric/percentage
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gold Price")

$ws.Range("B2").Value = "'48501.00"
$ws.Range("C2").Value = "'-402.00"
$ws.Range("D2").Value = "'-0.82%"

$ws.Range("B3").Value = "'48815.00"
$ws.Range("C3").Value = "'-410.00"
$ws.Range("D3").Value = "'-0.83%"
